$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 874.75
$ws.Range("I32").Value = 799.5
$ws.Range("J32").Value = 950
$ws.Range("K32").Value = 799.5
$ws.Range("L32").Value = 950
$ws.Range("M32").Value = -473.5
$ws.Range("N32").Value = -1602

$ws.Range("H111").Value = 1169
$ws.Range("I111").Value = 849.75
$ws.Range("J111").Value = 5000
$ws.Range("K111").Value = 2549.25
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = 517.75
$ws.Range("N111").Value = -21134

$ws.Range("H138").Value = 6316.7095
$ws.Range("I138").Value = 7962.2144
$ws.Range("J138").Value = 4961.5884
$ws.Range("K138").Value = 23886.6432
$ws.Range("L138").Value = 14884.7652
$ws.Range("M138").Value = -18746.6432
$ws.Range("N138").Value = -25164.7652

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14825.365
$ws.Range("I32").Value = 6077.614
$ws.Range("J32").Value = 27655.4
$ws.Range("K32").Value = 6077.614
$ws.Range("L32").Value = 27655.4
$ws.Range("M32").Value = -5790.614
$ws.Range("N32").Value = -28229.4

$ws.Range("H74").Value = 5020.1333
$ws.Range("I74").Value = 1477.4
$ws.Range("J74").Value = 6791.5
$ws.Range("K74").Value = 1477.4
$ws.Range("L74").Value = 6791.5
$ws.Range("M74").Value = -603.4000000000001
$ws.Range("N74").Value = -8539.5

$ws.Range("H77").Value = 5020.1333
$ws.Range("I77").Value = 1477.4
$ws.Range("J77").Value = 6791.5
$ws.Range("K77").Value = 7387
$ws.Range("L77").Value = 33957.5
$ws.Range("M77").Value = -3019
$ws.Range("N77").Value = -42693.5

$ws.Range("H97").Value = 812.5714
$ws.Range("I97").Value = 861.8
$ws.Range("J97").Value = 689.5
$ws.Range("K97").Value = 861.8
$ws.Range("L97").Value = 689.5
$ws.Range("M97").Value = -365.8
$ws.Range("N97").Value = -1681.5

$ws.Range("H102").Value = 1018.2759
$ws.Range("I102").Value = 986.2381
$ws.Range("J102").Value = 1102.375
$ws.Range("K102").Value = 986.2381
$ws.Range("L102").Value = 1102.375
$ws.Range("M102").Value = 635.7619
$ws.Range("N102").Value = -4346.375

$ws.Range("H132").Value = 3343.1316
$ws.Range("I132").Value = 1567.9333
$ws.Range("J132").Value = 10000.125
$ws.Range("K132").Value = 4703.7999
$ws.Range("L132").Value = 30000.375
$ws.Range("M132").Value = -2173.7999
$ws.Range("N132").Value = -35060.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1468.8
$ws.Range("I64").Value = 1299
$ws.Range("J64").Value = 1638.6
$ws.Range("K64").Value = 1299
$ws.Range("L64").Value = 1638.6
$ws.Range("M64").Value = -1074
$ws.Range("N64").Value = -2088.6

$ws.Range("H67").Value = 1468.8
$ws.Range("I67").Value = 1299
$ws.Range("J67").Value = 1638.6
$ws.Range("K67").Value = 1299
$ws.Range("L67").Value = 1638.6
$ws.Range("M67").Value = -519
$ws.Range("N67").Value = -3198.6

$ws.Range("H86").Value = 3576.5
$ws.Range("I86").Value = 2477.3333
$ws.Range("J86").Value = 4675.6665
$ws.Range("K86").Value = 2477.3333
$ws.Range("L86").Value = 4675.6665
$ws.Range("M86").Value = -1354.3333
$ws.Range("N86").Value = -6921.6665

$ws.Range("H89").Value = 3576.5
$ws.Range("I89").Value = 2477.3333
$ws.Range("J89").Value = 4675.6665
$ws.Range("K89").Value = 12386.6665
$ws.Range("L89").Value = 23378.3325
$ws.Range("M89").Value = -6770.666499999999
$ws.Range("N89").Value = -34610.3325

$ws.Range("H99").Value = 1370.4
$ws.Range("I99").Value = 1245.4445
$ws.Range("J99").Value = 2495
$ws.Range("K99").Value = 1245.4445
$ws.Range("L99").Value = 2495
$ws.Range("M99").Value = 252.5554999999999
$ws.Range("N99").Value = -5491

$ws.Range("H134").Value = 2066.0322
$ws.Range("I134").Value = 1539.2963
$ws.Range("J134").Value = 5621.5
$ws.Range("K134").Value = 4617.8889
$ws.Range("L134").Value = 16864.5
$ws.Range("M134").Value = -2082.8889
$ws.Range("N134").Value = -21934.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1999.5
$ws.Range("I2").Value = 2599
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 2599
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -2486
$ws.Range("N2").Value = -1626

$ws.Range("H94").Value = 1352.2
$ws.Range("I94").Value = 1022
$ws.Range("J94").Value = 1572.3334
$ws.Range("K94").Value = 1022
$ws.Range("L94").Value = 1572.3334
$ws.Range("M94").Value = -571
$ws.Range("N94").Value = -2474.3334

$ws.Range("H99").Value = 16081.615
$ws.Range("I99").Value = 14537.412
$ws.Range("J99").Value = 18998.445
$ws.Range("K99").Value = 14537.412
$ws.Range("L99").Value = 18998.445
$ws.Range("M99").Value = -13039.412
$ws.Range("N99").Value = -21994.445

$ws.Range("H126").Value = 16081.615
$ws.Range("I126").Value = 14537.412
$ws.Range("J126").Value = 18998.445
$ws.Range("K126").Value = 43612.236
$ws.Range("L126").Value = 56995.335
$ws.Range("M126").Value = -41142.236
$ws.Range("N126").Value = -61935.335

$ws.Range("H132").Value = 848.35297
$ws.Range("I132").Value = 848.35297
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2545.05891
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15.0589100000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2910984.5
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 4364977
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 13094931
$ws.Range("M32").Value = -8717
$ws.Range("N32").Value = -13095497

$ws.Range("H81").Value = 2627
$ws.Range("I81").Value = 1949.3334
$ws.Range("J81").Value = 3304.6667
$ws.Range("K81").Value = 5848.0002
$ws.Range("L81").Value = 9914.000100000001
$ws.Range("M81").Value = -4725.0002
$ws.Range("N81").Value = -12160.0001

$ws.Range("H84").Value = 2627
$ws.Range("I84").Value = 1949.3334
$ws.Range("J84").Value = 3304.6667
$ws.Range("K84").Value = 17544.0006
$ws.Range("L84").Value = 29742.0003
$ws.Range("M84").Value = -11928.0006
$ws.Range("N84").Value = -40974.0003

$ws.Range("H122").Value = 2220
$ws.Range("I122").Value = 2166.6667
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 19500.0003
$ws.Range("L122").Value = 20700
$ws.Range("M122").Value = -17050.0003
$ws.Range("N122").Value = -25600

$ws.Range("H129").Value = 1969.1666
$ws.Range("I129").Value = 2305
$ws.Range("J129").Value = 290
$ws.Range("K129").Value = 6915
$ws.Range("L129").Value = 870
$ws.Range("M129").Value = -1915
$ws.Range("N129").Value = -10870

$ws.Range("H131").Value = 1770.3636
$ws.Range("I131").Value = 746.5
$ws.Range("J131").Value = 2999
$ws.Range("K131").Value = 2239.5
$ws.Range("L131").Value = 8997
$ws.Range("M131").Value = 2800.5
$ws.Range("N131").Value = -19077

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8152.6
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 8440.75
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 8440.75
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -10436.75

$ws.Range("H83").Value = 8152.6
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 8440.75
$ws.Range("K83").Value = 35000
$ws.Range("L83").Value = 42203.75
$ws.Range("M83").Value = -30008
$ws.Range("N83").Value = -52187.75

$ws.Range("H97").Value = 1436.3462
$ws.Range("I97").Value = 1577.3684
$ws.Range("J97").Value = 1053.5714
$ws.Range("K97").Value = 1577.3684
$ws.Range("L97").Value = 1053.5714
$ws.Range("M97").Value = -1081.3684
$ws.Range("N97").Value = -2045.5714

$ws.Range("H107").Value = 65.8
$ws.Range("I107").Value = 75
$ws.Range("J107").Value = 52
$ws.Range("K107").Value = 75
$ws.Range("L107").Value = 52
$ws.Range("M107").Value = 1845
$ws.Range("N107").Value = -3892

$ws.Range("H113").Value = 4499.4546
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 4777.1113
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 4777.1113
$ws.Range("M113").Value = -1080
$ws.Range("N113").Value = -9117.1113

$ws.Range("H114").Value = 60000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 60000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws.Range("H132").Value = 2567.0625
$ws.Range("I132").Value = 2145.6667
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 6437.000100000001
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -3907.000100000001
$ws.Range("N132").Value = -31724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 30000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 30000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 30000
$ws.Range("N5").Value = -30226

$ws.Range("H22").Value = 1149.0625
$ws.Range("I22").Value = 703.1667
$ws.Range("J22").Value = 1416.6
$ws.Range("K22").Value = 703.1667
$ws.Range("L22").Value = 1416.6
$ws.Range("M22").Value = -408.1667
$ws.Range("N22").Value = -2006.6

$ws.Range("H27").Value = 1149.0625
$ws.Range("I27").Value = 703.1667
$ws.Range("J27").Value = 1416.6
$ws.Range("K27").Value = 703.1667
$ws.Range("L27").Value = 1416.6
$ws.Range("M27").Value = -596.1667
$ws.Range("N27").Value = -1630.6

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 2992.1538
$ws.Range("I46").Value = 1983.1666
$ws.Range("J46").Value = 3857
$ws.Range("K46").Value = 1983.1666
$ws.Range("L46").Value = 3857
$ws.Range("M46").Value = -1795.1666
$ws.Range("N46").Value = -4233

$ws.Range("H55").Value = 778.1875
$ws.Range("I55").Value = 332.375
$ws.Range("J55").Value = 1224
$ws.Range("K55").Value = 332.375
$ws.Range("L55").Value = 1224
$ws.Range("M55").Value = -159.375
$ws.Range("N55").Value = -1570

$ws.Range("H100").Value = 6820.2
$ws.Range("I100").Value = 2750.5
$ws.Range("J100").Value = 9533.333000000001
$ws.Range("K100").Value = 2750.5
$ws.Range("L100").Value = 9533.333000000001
$ws.Range("M100").Value = -2209.5
$ws.Range("N100").Value = -10615.333

$ws.Range("H136").Value = 3355.25
$ws.Range("I136").Value = 2807.3333
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 8421.999899999999
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -5871.999899999999
$ws.Range("N136").Value = -20097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5231.231
$ws.Range("I122").Value = 2240.7
$ws.Range("J122").Value = 15199.667
$ws.Range("K122").Value = 6722.099999999999
$ws.Range("L122").Value = 45599.001
$ws.Range("M122").Value = -4272.099999999999
$ws.Range("N122").Value = -50499.001

$ws.Range("H126").Value = 61705.707
$ws.Range("I126").Value = 168904.5
$ws.Range("J126").Value = 3233.6365
$ws.Range("K126").Value = 506713.5
$ws.Range("L126").Value = 9700.9095
$ws.Range("M126").Value = -504243.5
$ws.Range("N126").Value = -14640.9095
